$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.200.67'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.54'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7091'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.44'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3075'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07689'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.97'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08252'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.865.04'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7160'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.209'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.206.05'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.843'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.21'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007798'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.13'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.110.19'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.966'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1573'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +7.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.39'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.898'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.325'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.494'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.345'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.092'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05180'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.898'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.173'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7278'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.684'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.688'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.144.28'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.083'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.19'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.57'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5270'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.007.54'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.769'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000120'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.290'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.862'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.55%  '
